# Add "name" and "at" columns (name of the bus-type / line index) in front
# of the existing kW / kVAR columns, and dump results to an output folder
# (per commit message: "add names to lines & buses + dump the results to
# output folder"). This script only performs the spreadsheet-visible part
# of that change: inserting two new leading columns with the name/index
# data, shifting the existing kW/kVAR data two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:B (kW, kVAR) two columns to the right, to C:D,
# by inserting two new blank columns at the front.
$ws.Columns("A:B").Insert()

# New header row - column A first ("name"), column B ("at") filled in later
# (after the lorem/ipsum data) so new shared strings are registered in the
# same order as the authored workbook.
$ws.Range("A1").Value = "name"

# New data rows: alternating bus/line name in column A, sequential index
# in column B, for the 33 data rows that follow the header.
for ($i = 0; $i -lt 33; $i++) {
    $r = $i + 2
    if ($i % 2 -eq 0) {
        $name = "lorem"
    } else {
        $name = "ipsum"
    }
    $ws.Cells.Item($r, 1).Value = $name
}

$ws.Range("B1").Value = "at"

for ($i = 0; $i -lt 33; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $i + 1
}

# Selection, matching the saved workbook state
$ws.Range("R18").Select()
